# Apply the cryptos list refresh described by the diff.
# Updates price/volume figures and fixes three coin rows whose
# rank-order was swapped (WrappedEther/Polkadot, Filecoin/InternetComputer(DFINITY), Aptos/TrustWalletToken).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($Range, [string]$Value)
    # Force the cell to remain plain text even when the value looks like
    # a number (e.g. "1.006" or "54.69"), then restore the default "Normal"
    # style so no stray per-cell number format is left behind.
    $Range.NumberFormat = "@"
    $Range.Value = $Value
    $Range.Style = "Normal"
}

$ws.Range("D2").Value = '30.381.54'
$ws.Range("D3").Value = '2.104.79'
$ws.Range("E3").Value = '  +0.35%  '
$ws.Range("E4").Value = '  -0.16%  '
Set-TextCell $ws.Range("D5") '344.45'
$ws.Range("E5").Value = '  +0.37%  '
Set-TextCell $ws.Range("D6") '1.006'
$ws.Range("E6").Value = '  -0.12%  '
$ws.Range("E7").Value = '  +1.77%  '
Set-TextCell $ws.Range("D8") '0.4431'
$ws.Range("E8").Value = '  +0.70%  '
Set-TextCell $ws.Range("D9") '54.69'
$ws.Range("E9").Value = '  +2.87%  '
Set-TextCell $ws.Range("D10") '0.09379'
$ws.Range("E10").Value = '  +2.34%  '
Set-TextCell $ws.Range("D11") '1.172'
$ws.Range("E11").Value = '  +0.18%  '
Set-TextCell $ws.Range("D12") '24.92'
$ws.Range("E12").Value = '  +0.00%  '
Set-TextCell $ws.Range("D13") '8.663'
$ws.Range("E13").Value = '  +6.08%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '2.143.73'
$ws.Range("E14").Value = '  +2.04%  '
$ws.Range("B15").Value = 'Polkadot'
$ws.Range("C15").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextCell $ws.Range("D15") '6.927'
$ws.Range("E15").Value = '  +2.59%  '
Set-TextCell $ws.Range("D16") '101.78'
$ws.Range("E16").Value = '  +2.51%  '
$ws.Range("E17").Value = '  +1.13%  '
Set-TextCell $ws.Range("D18") '1.007'
$ws.Range("E18").Value = '  -0.14%  '
Set-TextCell $ws.Range("D19") '21.22'
$ws.Range("E19").Value = '  +1.87%  '
Set-TextCell $ws.Range("D20") '0.06715'
$ws.Range("E20").Value = '  +1.11%  '
Set-TextCell $ws.Range("D21") '6.385'
$ws.Range("E21").Value = '  +3.36%  '
$ws.Range("E22").Value = '  -0.14%  '
$ws.Range("D23").Value = '30.412.18'
$ws.Range("E23").Value = '  +2.23%  '
Set-TextCell $ws.Range("D24") '12.59'
$ws.Range("E24").Value = '  +0.27%  '
Set-TextCell $ws.Range("D25") '2.302'
$ws.Range("E25").Value = '  +0.55%  '
Set-TextCell $ws.Range("D26") '21.90'
$ws.Range("E26").Value = '  +0.53%  '
Set-TextCell $ws.Range("D27") '162.65'
$ws.Range("E27").Value = '  +0.04%  '
Set-TextCell $ws.Range("D28") '2.517'
$ws.Range("E28").Value = '  +0.35%  '
Set-TextCell $ws.Range("D29") '133.56'
$ws.Range("E29").Value = '  +0.91%  '
Set-TextCell $ws.Range("D30") '1.141'
$ws.Range("E30").Value = '  +0.77%  '
Set-TextCell $ws.Range("D31") '1.687'
$ws.Range("E31").Value = '  +3.19%  '
$ws.Range("E32").Value = '  +0.57%  '
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextCell $ws.Range("D33") '6.240'
$ws.Range("E33").Value = '  +1.46%  '
$ws.Range("B34").Value = 'InternetComputer(DFINITY)'
$ws.Range("C34").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextCell $ws.Range("D34") '6.721'
$ws.Range("E34").Value = '  +10.84%  '
Set-TextCell $ws.Range("D35") '3.920'
$ws.Range("E35").Value = '  -1.26%  '
Set-TextCell $ws.Range("D36") '10.31'
$ws.Range("E36").Value = '  +1.55%  '
Set-TextCell $ws.Range("D37") '0.02637'
$ws.Range("E37").Value = '  +2.97%  '
$ws.Range("E38").Value = '  +1.31%  '
Set-TextCell $ws.Range("D39") '0.7042'
$ws.Range("E39").Value = '  +2.87%  '
$ws.Range("B40").Value = 'Aptos'
$ws.Range("C40").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextCell $ws.Range("D40") '12.56'
$ws.Range("E40").Value = '  +1.44%  '
$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextCell $ws.Range("D41") '1.343'
$ws.Range("E41").Value = '  +3.46%  '
$ws.Range("E42").Value = '  +0.57%  '
Set-TextCell $ws.Range("D43") '0.6854'
$ws.Range("E43").Value = '  +3.16%  '
Set-TextCell $ws.Range("D44") '14.54'
$ws.Range("E44").Value = '  +2.59%  '
Set-TextCell $ws.Range("D45") '2.353'
$ws.Range("E45").Value = '  +2.56%  '
$ws.Range("E46").Value = '  -0.09%  '
Set-TextCell $ws.Range("D47") '1.400'
$ws.Range("E47").Value = '  +20.79%  '
Set-TextCell $ws.Range("D48") '3.644'
$ws.Range("E48").Value = '  +0.95%  '
Set-TextCell $ws.Range("D49") '0.00000000358'
$ws.Range("E49").Value = '  +4.97%  '
Set-TextCell $ws.Range("D50") '1.215'
$ws.Range("E50").Value = '  +9.41%  '
$ws.Range("E51").Value = '  +0.12%  '
